# LOM3050.xlsx update
#
# The sheet had several label/value pairs in columns B and C that were
# populated with the wrong text (each value was shifted up relative to
# its label). This change:
#   - inserts two rows (13-14) to hold the "Docentes responsaveis" names,
#   - re-populates every body-text cell in columns B/C from row 10 down
#     with the text that actually belongs next to its row-10..23 label,
#   - appends two brand-new rows (22-23) for "Norma de recuperação:" and
#     "Bibliografia:".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two rows at 13 to make room for the professors' names -----
$ws.Rows.Item(13).Resize(2).Insert()

# The inserted rows copied the formatting of the old row 13 ("Programa
# resumido:") onto column A; clear those A cells completely since the
# new rows 13/14 only hold data in columns B/C.
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()

# Give B13:C14 the same look (style) as the other data cells in columns
# B/C before putting values into them.
$ws.Range("B15").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B15").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C14").PasteSpecial(-4122)

# --- New rows 22/23: "Norma de recuperação:" and "Bibliografia:" ------
# Seed them (values + row height) from row 21 ("Bibliografia:" before the
# insert shifted things around) so they inherit the usual A/B/C styling.
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C21").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("A21").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("B21").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("C21").Copy()
$ws.Range("C23").PasteSpecial(-4122)

$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120

# --- Now (re)populate every body-text cell in B/C for rows 10,13-23 ---

$objetivosText = "Complementar a formação dos alunos em Engenharia de Materiais abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."
$ws.Range("B10").Value = $objetivosText
$ws.Range("C10").Value = $objetivosText

$ws.Range("B13").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C13").Value = "471420 - Carlos Antonio Reis Pereira Baptista"

$ws.Range("B14").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C14").Value = "519033 - Carlos Yujiro Shigue"

$programaResumidoText = "A definir, de acordo com o tópico programado."
$ws.Range("B15").Value = $programaResumidoText
$ws.Range("C15").Value = $programaResumidoText

$programaText = "O conteúdo desta disciplina (optativa)será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."
$ws.Range("B17").Value = $programaText
$ws.Range("C17").Value = $programaText

$metodoText = "Este curso deverá conter avaliações escritas e desenvolvimento de Estudo de Casos ou Projetos na área de Engenharia de Materiais. Sendo necessário aplicar pelo menos dois tipos de avaliações diferentes."
$ws.Range("B20").Value = $metodoText
$ws.Range("C20").Value = $metodoText

$criterioText = "A média do semestre será computada com base na relação: M=(A1+A2)/2"
$ws.Range("B21").Value = $criterioText
$ws.Range("C21").Value = $criterioText

$ws.Range("A22").Value = "Norma de recuperação:"
$normaText = "Não cabe recuperação."
$ws.Range("B22").Value = $normaText
$ws.Range("C22").Value = $normaText

$ws.Range("A23").Value = "Bibliografia:"
$bibliografiaText = "Apostila ou texto fornecido pelo(s) docente(s) responsáveis.Artigos extraídos de revistas especializadas na área de Ciência e Engenharia de Materiais."
$ws.Range("B23").Value = $bibliografiaText
$ws.Range("C23").Value = $bibliografiaText
